# Applies the crypto price/volume refresh described in the commit diff.
# Each row 2-51 holds one coin: B=Coin, C=Link, D=Price, E=Volume(1h).
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel keeps them as text (matching the original inlineStr cells) instead
# of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '50.865.85'
$ws.Cells.Item(2, 5).Value = '  -2.65%  '

$ws.Cells.Item(3, 4).Value = '2.881.08'
$ws.Cells.Item(3, 5).Value = '  -2.75%  '

$ws.Cells.Item(4, 5).Value = '  -0.08%  '

$ws.Cells.Item(5, 4).Value = '''365.79'
$ws.Cells.Item(5, 5).Value = '  +2.66%  '

$ws.Cells.Item(6, 4).Value = '''101.41'
$ws.Cells.Item(6, 5).Value = '  -7.13%  '

$ws.Cells.Item(7, 4).Value = '''0.535'
$ws.Cells.Item(7, 5).Value = '  -5.80%  '

$ws.Cells.Item(8, 5).Value = '  -0.04%  '

$ws.Cells.Item(9, 4).Value = '''0.584'
$ws.Cells.Item(9, 5).Value = '  -7.01%  '

$ws.Cells.Item(10, 4).Value = '''36.42'
$ws.Cells.Item(10, 5).Value = '  -6.29%  '

$ws.Cells.Item(11, 5).Value = '  +0.64%  '

$ws.Cells.Item(12, 4).Value = '''0.0829'
$ws.Cells.Item(12, 5).Value = '  -5.01%  '

$ws.Cells.Item(13, 4).Value = '''18.20'
$ws.Cells.Item(13, 5).Value = '  -5.83%  '

$ws.Cells.Item(14, 4).Value = '3.338.16'
$ws.Cells.Item(14, 5).Value = '  -3.01%  '

$ws.Cells.Item(15, 4).Value = '''7.35'
$ws.Cells.Item(15, 5).Value = '  -5.33%  '

$ws.Cells.Item(16, 4).Value = '2.876.31'
$ws.Cells.Item(16, 5).Value = '  -2.60%  '

$ws.Cells.Item(17, 4).Value = '''0.925'
$ws.Cells.Item(17, 5).Value = '  -5.77%  '

$ws.Cells.Item(18, 4).Value = '50.810.84'
$ws.Cells.Item(18, 5).Value = '  -2.77%  '

$ws.Cells.Item(19, 5).Value = '  -7.51%  '

$ws.Cells.Item(20, 4).Value = '''7.14'
$ws.Cells.Item(20, 5).Value = '  -5.87%  '

$ws.Cells.Item(21, 4).Value = '''12.75'
$ws.Cells.Item(21, 5).Value = '  -7.50%  '

$ws.Cells.Item(22, 4).Value = '0.0₃0937'
$ws.Cells.Item(22, 5).Value = '  -4.68%  '

$ws.Cells.Item(23, 4).Value = '''67.71'
$ws.Cells.Item(23, 5).Value = '  -3.67%  '

$ws.Cells.Item(24, 4).Value = '''256.72'
$ws.Cells.Item(24, 5).Value = '  -4.32%  '

$ws.Cells.Item(25, 5).Value = '  -4.67%  '

$ws.Cells.Item(26, 2).Value = 'Kaspa'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(26, 4).Value = '''0.169'
$ws.Cells.Item(26, 5).Value = '  -4.74%  '

$ws.Cells.Item(27, 2).Value = 'Dai'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(27, 4).Value = '''1.00'
$ws.Cells.Item(27, 5).Value = '  +0.12%  '

$ws.Cells.Item(28, 4).Value = '''25.45'
$ws.Cells.Item(28, 5).Value = '  -6.18%  '

$ws.Cells.Item(29, 4).Value = '''6.89'
$ws.Cells.Item(29, 5).Value = '  -9.87%  '

$ws.Cells.Item(30, 5).Value = '  -5.36%  '

$ws.Cells.Item(31, 5).Value = '  -5.87%  '

$ws.Cells.Item(32, 5).Value = '  -2.31%  '

$ws.Cells.Item(33, 5).Value = '  -2.79%  '

$ws.Cells.Item(34, 4).Value = '''34.24'
$ws.Cells.Item(34, 5).Value = '  -8.04%  '

$ws.Cells.Item(35, 4).Value = '''50.58'
$ws.Cells.Item(35, 5).Value = '  -3.53%  '

$ws.Cells.Item(36, 5).Value = '  +0.23%  '

$ws.Cells.Item(37, 5).Value = '  -6.41%  '

$ws.Cells.Item(38, 4).Value = '''3.02'
$ws.Cells.Item(38, 5).Value = '  -5.81%  '

$ws.Cells.Item(39, 4).Value = '''2.62'
$ws.Cells.Item(39, 5).Value = '  -3.39%  '

$ws.Cells.Item(40, 4).Value = '''16.81'
$ws.Cells.Item(40, 5).Value = '  -6.79%  '

$ws.Cells.Item(41, 5).Value = '  -9.71%  '

$ws.Cells.Item(42, 5).Value = '  -6.10%  '

$ws.Cells.Item(43, 4).Value = '''21.77'
$ws.Cells.Item(43, 5).Value = '  -7.29%  '

$ws.Cells.Item(44, 4).Value = '''117.64'
$ws.Cells.Item(44, 5).Value = '  -1.28%  '

$ws.Cells.Item(45, 5).Value = '  -3.69%  '

$ws.Cells.Item(46, 4).Value = '2.016.86'
$ws.Cells.Item(46, 5).Value = '  -5.64%  '

$ws.Cells.Item(47, 5).Value = '  -6.23%  '

$ws.Cells.Item(48, 4).Value = '''3.13'
$ws.Cells.Item(48, 5).Value = '  -8.77%  '

$ws.Cells.Item(49, 4).Value = '3.179.32'
$ws.Cells.Item(49, 5).Value = '  -2.68%  '

$ws.Cells.Item(50, 5).Value = '  -3.90%  '

$ws.Cells.Item(51, 4).Value = '''0.0307'
$ws.Cells.Item(51, 5).Value = '  -12.24%  '
